$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2053571428571428
$ws.Range("C2").Value = 0.5357142857142857
$ws.Range("J2").Value = 0.02380952380952381
$ws.Range("P2").Value = 0.1488095238095238
$ws.Range("S2").Value = 0.08630952380952381
$ws.Range("B3").Value = 0.02061855670103093
$ws.Range("C3").Value = 0.07216494845360824
$ws.Range("J3").Value = 0.02577319587628866
$ws.Range("P3").Value = 0.7010309278350515
$ws.Range("S3").Value = 0.1804123711340206
$ws.Range("J4").Value = 0.1020408163265306
$ws.Range("P4").Value = 0.6122448979591837
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("P5").Value = 0.5714285714285714
$ws.Range("S5").Value = 0.4285714285714285
$ws.Range("B6").Value = 0.06557377049180328
$ws.Range("D6").Value = 0.004098360655737705
$ws.Range("F6").Value = 0.04508196721311476
$ws.Range("J6").Value = 0.2418032786885246
$ws.Range("O6").Value = 0.01229508196721311
$ws.Range("Q6").Value = 0.1065573770491803
$ws.Range("R6").Value = 0.09016393442622951
$ws.Range("S6").Value = 0.4344262295081967
$ws.Range("B7").Value = 0.07657657657657657
$ws.Range("D7").Value = 0.02702702702702703
$ws.Range("F7").Value = 0.06756756756756757
$ws.Range("J7").Value = 0.0990990990990991
$ws.Range("O7").Value = 0.01801801801801802
$ws.Range("Q7").Value = 0.1891891891891892
$ws.Range("R7").Value = 0.0990990990990991
$ws.Range("S7").Value = 0.4234234234234234
$ws.Range("B8").Value = 0.1048543689320388
$ws.Range("D8").Value = 0.02330097087378641
$ws.Range("E8").Value = 0.005825242718446602
$ws.Range("F8").Value = 0.05631067961165048
$ws.Range("J8").Value = 0.1339805825242718
$ws.Range("O8").Value = 0.02135922330097087
$ws.Range("Q8").Value = 0.1533980582524272
$ws.Range("R8").Value = 0.1067961165048544
$ws.Range("S8").Value = 0.3941747572815534
$ws.Range("B9").Value = 0.09767441860465116
$ws.Range("D9").Value = 0.02790697674418605
$ws.Range("F9").Value = 0.08372093023255814
$ws.Range("J9").Value = 0.1162790697674419
$ws.Range("O9").Value = 0.0186046511627907
$ws.Range("Q9").Value = 0.1627906976744186
$ws.Range("R9").Value = 0.1069767441860465
$ws.Range("S9").Value = 0.386046511627907
$ws.Range("B10").Value = 0.1130820399113082
$ws.Range("D10").Value = 0.01847745750184775
$ws.Range("E10").Value = 0.002956393200295639
$ws.Range("F10").Value = 0.07760532150776053
$ws.Range("J10").Value = 0.1160384331116038
$ws.Range("O10").Value = 0.01404286770140429
$ws.Range("Q10").Value = 0.1677753141167775
$ws.Range("R10").Value = 0.1027346637102735
$ws.Range("S10").Value = 0.3872875092387287
$ws.Range("F11").Value = 0.002932551319648094
$ws.Range("G11").Value = 0.1319648093841642
$ws.Range("J11").Value = 0.07624633431085044
$ws.Range("K11").Value = 0.187683284457478
$ws.Range("L11").Value = 0.5835777126099707
$ws.Range("S11").Value = 0.01759530791788856
$ws.Range("G12").Value = 0.7403846153846154
$ws.Range("J12").Value = 0.2067307692307692
$ws.Range("L12").Value = 0.02403846153846154
$ws.Range("S12").Value = 0.02884615384615385
$ws.Range("G13").Value = 0.6875
$ws.Range("J13").Value = 0.2291666666666667
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.01260504201680672
$ws.Range("H15").Value = 0.2310924369747899
$ws.Range("I15").Value = 0.04621848739495799
$ws.Range("J15").Value = 0.3319327731092437
$ws.Range("K15").Value = 0.0546218487394958
$ws.Range("M15").Value = 0.004201680672268907
$ws.Range("O15").Value = 0.04201680672268908
$ws.Range("S15").Value = 0.2773109243697479
$ws.Range("F16").Value = 0.009433962264150943
$ws.Range("H16").Value = 0.1839622641509434
$ws.Range("I16").Value = 0.07075471698113207
$ws.Range("J16").Value = 0.4056603773584906
$ws.Range("K16").Value = 0.1509433962264151
$ws.Range("M16").Value = 0.009433962264150943
$ws.Range("O16").Value = 0.05660377358490566
$ws.Range("S16").Value = 0.1132075471698113
$ws.Range("F17").Value = 0.02205882352941177
$ws.Range("H17").Value = 0.2009803921568628
$ws.Range("I17").Value = 0.1102941176470588
$ws.Range("J17").Value = 0.3725490196078431
$ws.Range("K17").Value = 0.1127450980392157
$ws.Range("M17").Value = 0.0196078431372549
$ws.Range("O17").Value = 0.04411764705882353
$ws.Range("S17").Value = 0.1176470588235294
$ws.Range("F18").Value = 0.03065134099616858
$ws.Range("H18").Value = 0.2068965517241379
$ws.Range("I18").Value = 0.06513409961685823
$ws.Range("J18").Value = 0.4022988505747127
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("M18").Value = 0.02298850574712644
$ws.Range("N18").Value = 0.003831417624521073
$ws.Range("O18").Value = 0.06130268199233716
$ws.Range("S18").Value = 0.09578544061302682
$ws.Range("F19").Value = 0.01275690999291283
$ws.Range("H19").Value = 0.2034018426647768
$ws.Range("I19").Value = 0.09071580439404678
$ws.Range("J19").Value = 0.371367824238129
$ws.Range("K19").Value = 0.109851169383416
$ws.Range("M19").Value = 0.0219702338766832
$ws.Range("N19").Value = 0.0007087172218284905
$ws.Range("O19").Value = 0.07725017717930546
$ws.Range("S19").Value = 0.1119773210489015
